# Fix the title on the first slide ("P1"): "Part III: Rechtliche Grundlagen"
# should read "Part I: Rechtliche Grundlagen". The original text is a single
# run; we split it into "Part " / "I: " / "Rechtliche Grundlagen" by
# replacing the "III: " substring with "I: " via a Characters() sub-range,
# which is how PowerPoint naturally breaks a run into multiple runs at an
# edit boundary.

$p = $ppt.ActivePresentation

$needle = "Part III: Rechtliche Grundlagen"
$replacement = "I: "   # replaces "III: " (chars 6-10 of $needle)

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }
        if (-not $shape.TextFrame.HasText) { continue }

        $tr = $shape.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf($needle)
        if ($idx -ge 0) {
            # 1-indexed start of "Part III: Rechtliche Grundlagen" within the TextRange
            $start = $idx + 1
            # "III: " starts 5 chars into the needle ("Part ") and is 5 chars long
            $replaceStart = $start + 5
            $old = $tr.Characters($replaceStart, 5)
            $old.Text = $replacement
        }
    }
}
